$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($Worksheet, [string]$Address, [string]$Text)
    $cell = $Worksheet.Range($Address)
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "42.777.78"
$ws.Range("E2").Value = "  +0.14%  "
$ws.Range("D3").Value = "2.565.68"
$ws.Range("E3").Value = "  +0.84%  "
Set-TextCell $ws "D4" "0.998"
$ws.Range("E4").Value = "  -0.08%  "
Set-TextCell $ws "D5" "302.88"
$ws.Range("E5").Value = "  +2.29%  "
Set-TextCell $ws "D6" "97.02"
$ws.Range("E6").Value = "  +3.89%  "
Set-TextCell $ws "D7" "0.574"
$ws.Range("E7").Value = "  +0.48%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  -0.07%  "
Set-TextCell $ws "D10" "36.11"
$ws.Range("E10").Value = "  +1.46%  "
Set-TextCell $ws "D11" "0.0808"
$ws.Range("E11").Value = "  +0.44%  "
$ws.Range("E12").Value = "  +10.17%  "
Set-TextCell $ws "D13" "7.54"
$ws.Range("E13").Value = "  -1.04%  "
$ws.Range("D14").Value = "2.633.96"
$ws.Range("E14").Value = "  +4.01%  "
Set-TextCell $ws "D15" "0.876"
$ws.Range("E15").Value = "  +1.30%  "
Set-TextCell $ws "D16" "14.40"
$ws.Range("E16").Value = "  +2.15%  "
$ws.Range("D17").Value = "42.831.12"
$ws.Range("E17").Value = "  +0.18%  "
Set-TextCell $ws "D18" "13.31"
$ws.Range("E18").Value = "  +6.79%  "
$ws.Range("D19").Value = "0.0₃0987"
$ws.Range("E19").Value = "  +2.36%  "
Set-TextCell $ws "D20" "6.62"
$ws.Range("E20").Value = "  +1.01%  "
Set-TextCell $ws "D21" "71.48"
$ws.Range("E21").Value = "  -0.75%  "
Set-TextCell $ws "D22" "257.12"
$ws.Range("E22").Value = "  -0.75%  "
Set-TextCell $ws "D23" "2.96"
$ws.Range("E23").Value = "  +2.83%  "
$ws.Range("E24").Value = "  -1.35%  "
Set-TextCell $ws "D25" "28.21"
$ws.Range("E25").Value = "  -4.25%  "
$ws.Range("E26").Value = "  -0.04%  "
Set-TextCell $ws "D27" "39.37"
$ws.Range("E27").Value = "  +9.42%  "
Set-TextCell $ws "D28" "10.08"
$ws.Range("E28").Value = "  +1.70%  "
$ws.Range("E29").Value = "  -1.56%  "
Set-TextCell $ws "D30" "6.00"
$ws.Range("E30").Value = "  +2.04%  "
Set-TextCell $ws "D31" "156.55"
$ws.Range("E31").Value = "  +4.33%  "
$ws.Range("E32").Value = "  +1.77%  "
Set-TextCell $ws "D33" "2.15"
$ws.Range("E33").Value = "  -0.23%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextCell $ws "D34" "0.0801"
$ws.Range("E34").Value = "  +1.02%  "
$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextCell $ws "D35" "3.31"
$ws.Range("E35").Value = "  -2.77%  "
Set-TextCell $ws "D36" "26.45"
$ws.Range("E36").Value = "  +8.88%  "
Set-TextCell $ws "D37" "0.115"
$ws.Range("E37").Value = "  +1.08%  "
Set-TextCell $ws "D38" "18.04"
$ws.Range("E38").Value = "  +13.64%  "
$ws.Range("E39").Value = "  +0.28%  "
Set-TextCell $ws "D40" "3.86"
$ws.Range("E40").Value = "  +1.84%  "
$ws.Range("E41").Value = "  +30.47%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextCell $ws "D42" "0.0307"
$ws.Range("E42").Value = "  -0.27%  "
$ws.Range("B43").Value = "NEARProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextCell $ws "D43" "3.37"
$ws.Range("E43").Value = "  -1.71%  "
$ws.Range("D44").Value = "2.065.93"
$ws.Range("E44").Value = "  -0.67%  "
$ws.Range("E45").Value = "  -0.01%  "
Set-TextCell $ws "D46" "88.30"
$ws.Range("E46").Value = "  +4.25%  "
Set-TextCell $ws "D47" "9.24"
$ws.Range("E47").Value = "  +5.60%  "
Set-TextCell $ws "D48" "76.58"
$ws.Range("E48").Value = "  +10.47%  "
$ws.Range("D49").Value = "2.810.33"
$ws.Range("E49").Value = "  +0.84%  "
Set-TextCell $ws "D50" "104.12"
$ws.Range("E50").Value = "  +1.15%  "
Set-TextCell $ws "D51" "0.189"
$ws.Range("E51").Value = "  +2.44%  "

Write-Output "Updated cryptos list"
